# Update cryptocurrency price/volume data per latest scrape
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.705.68'
$ws.Range('E2').Value = '  -1.71%  '
$ws.Range('D3').Value = '1.803.18'
$ws.Range('E3').Value = '  -1.26%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').Value = '''231.56'
$ws.Range('E5').Value = '  -1.79%  '
$ws.Range('D6').Value = '''0.5950'
$ws.Range('E6').Value = '  -2.52%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').Value = '''0.2783'
$ws.Range('E8').Value = '  -0.80%  '
$ws.Range('D9').Value = '''0.06838'
$ws.Range('E9').Value = '  -3.56%  '
$ws.Range('D10').Value = '''23.32'
$ws.Range('E10').Value = '  -0.83%  '
$ws.Range('D11').Value = '''0.07538'
$ws.Range('E11').Value = '  -1.63%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.806.74'
$ws.Range('E12').Value = '  -1.02%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').Value = '''4.774'
$ws.Range('E13').Value = '  -0.76%  '
$ws.Range('D14').Value = '''0.6240'
$ws.Range('E14').Value = '  -1.18%  '
$ws.Range('D15').Value = '2.048.75'
$ws.Range('E15').Value = '  -1.20%  '
$ws.Range('D16').Value = '''0.000009341'
$ws.Range('E16').Value = '  -6.58%  '
$ws.Range('D17').Value = '''75.55'
$ws.Range('E17').Value = '  -3.95%  '
$ws.Range('D18').Value = '28.662.98'
$ws.Range('E18').Value = '  -1.74%  '
$ws.Range('D19').Value = '''5.488'
$ws.Range('E19').Value = '  -6.38%  '
$ws.Range('E20').Value = '  +0.12%  '
$ws.Range('D21').Value = '''210.08'
$ws.Range('E21').Value = '  -7.46%  '
$ws.Range('D22').Value = '''11.45'
$ws.Range('E22').Value = '  -2.67%  '
$ws.Range('D23').Value = '''6.860'
$ws.Range('E23').Value = '  -2.01%  '
$ws.Range('E24').Value = '  +0.18%  '
$ws.Range('D25').Value = '''154.27'
$ws.Range('E25').Value = '  -0.81%  '
$ws.Range('D26').Value = '''7.866'
$ws.Range('E26').Value = '  -2.19%  '
$ws.Range('E27').Value = '  -2.80%  '
$ws.Range('D28').Value = '''16.38'
$ws.Range('E28').Value = '  -1.32%  '
$ws.Range('D29').Value = '''1.431'
$ws.Range('E29').Value = '  -4.15%  '
$ws.Range('D30').Value = '''0.06212'
$ws.Range('E30').Value = '  -2.19%  '
$ws.Range('E31').Value = '  -1.84%  '
$ws.Range('D32').Value = '''3.780'
$ws.Range('E32').Value = '  -1.06%  '
$ws.Range('D33').Value = '''3.750'
$ws.Range('E33').Value = '  -1.18%  '
$ws.Range('D34').Value = '''1.721'
$ws.Range('E34').Value = '  -1.19%  '
$ws.Range('E35').Value = '  -5.43%  '
$ws.Range('D36').Value = '''0.6395'
$ws.Range('E36').Value = '  -0.98%  '
$ws.Range('D37').Value = '''2.492'
$ws.Range('E37').Value = '  -2.10%  '
$ws.Range('E38').Value = '  -0.93%  '
$ws.Range('D39').Value = '''0.01716'
$ws.Range('E39').Value = '  -1.20%  '
$ws.Range('D40').Value = '''6.457'
$ws.Range('E40').Value = '  -1.37%  '
$ws.Range('D41').Value = '1.132.38'
$ws.Range('E41').Value = '  -6.63%  '
$ws.Range('D42').Value = '''0.8743'
$ws.Range('E42').Value = '  -4.32%  '
$ws.Range('E43').Value = '  +0.29%  '
$ws.Range('D44').Value = '''100.81'
$ws.Range('E44').Value = '  +0.11%  '
$ws.Range('D45').Value = '1.965.60'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('D46').Value = '''60.55'
$ws.Range('E46').Value = '  -3.34%  '
$ws.Range('E47').Value = '  -2.41%  '
$ws.Range('D48').Value = '''1.599'
$ws.Range('E48').Value = '  -0.31%  '
$ws.Range('D49').Value = '''0.05476'
$ws.Range('E49').Value = '  -0.83%  '
$ws.Range('D50').Value = '''8.333'
$ws.Range('E50').Value = '  -2.49%  '
$ws.Range('D51').Value = '''0.4484'
$ws.Range('E51').Value = '  -1.83%  '
